$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 67.77251700000001
$ws.Range("H2").Value = 203.317551
$ws.Range("I2").Value = 0.4079637943863715
$ws.Range("J2").Value = 0.4079637943863715
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.08181366666667
$ws.Range("N2").Value = 63.245441
$ws.Range("O2").Value = 0.0571606014598545
$ws.Range("P2").Value = 0.0571606014598545
$ws.Range("Q2").Value = 1428.767575114999
$ws.Range("R2").Value = 12858.90817603499
$ws.Range("S2").Value = 0.02331945586096941
$ws.Range("T2").Value = 0.02331945586096941
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 67.77251700000001
$ws.Range("H3").Value = 203.317551
$ws.Range("I3").Value = 0.4079637943863715
$ws.Range("J3").Value = 0.4079637943863715
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 301.6001486666667
$ws.Range("N3").Value = 904.800446
$ws.Range("O3").Value = 0.8177496571571792
$ws.Range("P3").Value = 0.8177496571571792
$ws.Range("Q3").Value = 20440.2012027142
$ws.Range("R3").Value = 183961.8108244278
$ws.Range("S3").Value = 0.3336122529919972
$ws.Range("T3").Value = 0.3336122529919972
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 67.77251700000001
$ws.Range("H4").Value = 203.317551
$ws.Range("I4").Value = 0.4079637943863715
$ws.Range("J4").Value = 0.4079637943863715
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 46.13524966666667
$ws.Range("N4").Value = 138.405749
$ws.Range("O4").Value = 0.1250897413829664
$ws.Range("P4").Value = 0.1250897413829664
$ws.Range("Q4").Value = 3126.701992333412
$ws.Range("R4").Value = 28140.3179310007
$ws.Range("S4").Value = 0.05103208553340487
$ws.Range("T4").Value = 0.05103208553340487
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 60.97760633333333
$ws.Range("H5").Value = 182.932819
$ws.Range("I5").Value = 0.3670611149405164
$ws.Range("J5").Value = 0.3670611149405164
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.08181366666667
$ws.Range("N5").Value = 63.245441
$ws.Range("O5").Value = 0.0571606014598545
$ws.Range("P5").Value = 0.0571606014598545
$ws.Range("Q5").Value = 1285.518534558686
$ws.Range("R5").Value = 11569.66681102818
$ws.Range("S5").Value = 0.0209814341025247
$ws.Range("T5").Value = 0.0209814341025247
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 60.97760633333333
$ws.Range("H6").Value = 182.932819
$ws.Range("I6").Value = 0.3670611149405164
$ws.Range("J6").Value = 0.3670611149405164
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 301.6001486666667
$ws.Range("N6").Value = 904.800446
$ws.Range("O6").Value = 0.8177496571571792
$ws.Range("P6").Value = 0.8177496571571792
$ws.Range("Q6").Value = 18390.85513547081
$ws.Range("R6").Value = 165517.6962192373
$ws.Range("S6").Value = 0.3001641008983392
$ws.Range("T6").Value = 0.3001641008983392
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 60.97760633333333
$ws.Range("H7").Value = 182.932819
$ws.Range("I7").Value = 0.3670611149405164
$ws.Range("J7").Value = 0.3670611149405164
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.13524966666667
$ws.Range("N7").Value = 138.405749
$ws.Range("O7").Value = 0.1250897413829664
$ws.Range("P7").Value = 0.1250897413829664
$ws.Range("Q7").Value = 2813.217092264048
$ws.Range("R7").Value = 25318.95383037643
$ws.Range("S7").Value = 0.04591557993965249
$ws.Range("T7").Value = 0.04591557993965248
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 37.37372866666667
$ws.Range("H8").Value = 112.121186
$ws.Range("I8").Value = 0.2249750906731122
$ws.Range("J8").Value = 0.2249750906731122
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.08181366666667
$ws.Range("N8").Value = 63.245441
$ws.Range("O8").Value = 0.0571606014598545
$ws.Range("P8").Value = 0.0571606014598545
$ws.Range("Q8").Value = 787.9059837792252
$ws.Range("R8").Value = 7091.153854013027
$ws.Range("S8").Value = 0.01285971149636039
$ws.Range("T8").Value = 0.01285971149636039
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 37.37372866666667
$ws.Range("H9").Value = 112.121186
$ws.Range("I9").Value = 0.2249750906731122
$ws.Range("J9").Value = 0.2249750906731122
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 301.6001486666667
$ws.Range("N9").Value = 904.800446
$ws.Range("O9").Value = 0.8177496571571792
$ws.Range("P9").Value = 0.8177496571571792
$ws.Range("Q9").Value = 11271.92212209433
$ws.Range("R9").Value = 101447.299098849
$ws.Range("S9").Value = 0.1839733032668428
$ws.Range("T9").Value = 0.1839733032668428
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 37.37372866666667
$ws.Range("H10").Value = 112.121186
$ws.Range("I10").Value = 0.2249750906731122
$ws.Range("J10").Value = 0.2249750906731122
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 46.13524966666667
$ws.Range("N10").Value = 138.405749
$ws.Range("O10").Value = 0.1250897413829664
$ws.Range("P10").Value = 0.1250897413829664
$ws.Range("Q10").Value = 1724.246303010924
$ws.Range("R10").Value = 15518.21672709832
$ws.Range("S10").Value = 0.02814207590990901
$ws.Range("T10").Value = 0.02814207590990901
